$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text number format on Price cells whose new values would otherwise
# be auto-parsed by Excel as numbers (e.g. "1.00", "211.95") so they are
# stored as text, matching the source data which is all inline text.
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D51").NumberFormat = "@"

$ws.Range("D2").Value = "26.775.88"
$ws.Range("E2").Value = "  +0.25%  "
$ws.Range("D3").Value = "1.605.17"
$ws.Range("E3").Value = "  +0.47%  "
$ws.Range("D4").Value = "1.00"
$ws.Range("E4").Value = "  +0.06%  "
$ws.Range("D5").Value = "211.95"
$ws.Range("E5").Value = "  +0.28%  "
$ws.Range("D6").Value = "0.513"
$ws.Range("E6").Value = "  -0.10%  "
$ws.Range("D7").Value = "1.00"
$ws.Range("E7").Value = "  +0.00%  "
$ws.Range("D8").Value = "0.0620"
$ws.Range("E8").Value = "  -0.01%  "
$ws.Range("E9").Value = "  -0.07%  "
$ws.Range("D10").Value = "19.72"
$ws.Range("E10").Value = "  +1.03%  "
$ws.Range("D11").Value = "0.0849"
$ws.Range("E11").Value = "  +0.90%  "
$ws.Range("D12").Value = "1.826.00"
$ws.Range("E12").Value = "  +0.21%  "
$ws.Range("D13").Value = "1.586.49"
$ws.Range("E13").Value = "  -0.64%  "
$ws.Range("E14").Value = "  +0.64%  "
$ws.Range("D15").Value = "0.527"
$ws.Range("E15").Value = "  +0.68%  "
$ws.Range("D16").Value = "65.24"
$ws.Range("E16").Value = "  -0.21%  "
$ws.Range("D17").Value = "0.0₃0743"
$ws.Range("E17").Value = "  -3.41%  "
$ws.Range("B18").Value = "BitcoinCash"
$ws.Range("C18").Value = "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"
$ws.Range("D18").Value = "209.92"
$ws.Range("E18").Value = "  +0.13%  "
$ws.Range("D19").Value = "7.19"
$ws.Range("E19").Value = "  +0.32%  "
$ws.Range("B20").Value = "Dai"
$ws.Range("C20").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D20").Value = "1.00"
$ws.Range("E20").Value = "  +0.04%  "
$ws.Range("E21").Value = "  +0.79%  "
$ws.Range("E22").Value = "  -2.98%  "
$ws.Range("D23").Value = "9.05"
$ws.Range("E23").Value = "  +1.15%  "
$ws.Range("D24").Value = "143.87"
$ws.Range("E24").Value = "  +0.51%  "
$ws.Range("D25").Value = "1.00"
$ws.Range("E25").Value = "  +0.12%  "
$ws.Range("D26").Value = "7.15"
$ws.Range("E26").Value = "  +0.10%  "
$ws.Range("E27").Value = "  -0.61%  "
$ws.Range("D28").Value = "15.40"
$ws.Range("E28").Value = "  +0.52%  "
$ws.Range("E29").Value = "  -1.72%  "
$ws.Range("E30").Value = "  -0.19%  "
$ws.Range("E31").Value = "  +0.60%  "
$ws.Range("D32").Value = "3.00"
$ws.Range("E32").Value = "  +1.02%  "
$ws.Range("E33").Value = "  +19.77%  "
$ws.Range("D34").Value = "1.281.81"
$ws.Range("E34").Value = "  -0.45%  "
$ws.Range("D35").Value = "2.49"
$ws.Range("E35").Value = "  +1.12%  "
$ws.Range("E36").Value = "  +0.51%  "
$ws.Range("D37").Value = "0.593"
$ws.Range("E37").Value = "  -4.07%  "
$ws.Range("E38").Value = "  -1.71%  "
$ws.Range("D39").Value = "0.828"
$ws.Range("E39").Value = "  +0.05%  "
$ws.Range("B40").Value = "MXToken"
$ws.Range("C40").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("D40").Value = "2.25"
$ws.Range("E40").Value = "  +3.10%  "
$ws.Range("B41").Value = "FraxShare"
$ws.Range("C41").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D41").Value = "5.47"
$ws.Range("E41").Value = "  +0.55%  "
$ws.Range("D42").Value = "0.779"
$ws.Range("E42").Value = "  -0.76%  "
$ws.Range("D43").Value = "62.95"
$ws.Range("E43").Value = "  -0.49%  "
$ws.Range("D44").Value = "1.739.19"
$ws.Range("E44").Value = "  +0.22%  "
$ws.Range("D45").Value = "90.51"
$ws.Range("E45").Value = "  -0.93%  "
$ws.Range("E46").Value = "  +0.43%  "
$ws.Range("D47").Value = "0.103"
$ws.Range("E47").Value = "  +2.16%  "
$ws.Range("E48").Value = "  +0.89%  "
$ws.Range("B49").Value = "BabyDogeCoin"
$ws.Range("C49").Value = "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
$ws.Range("D49").Value = "0.0₇0974"
$ws.Range("E49").Value = "  -8.25%  "
$ws.Range("B50").Value = "EnergySwap"
$ws.Range("C50").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D50").Value = "7.55"
$ws.Range("E50").Value = "  +3.12%  "
$ws.Range("B51").Value = "USDD"
$ws.Range("C51").Value = "https://coinranking.com/coin/z2PZIKQL7+usdd-usdd"
$ws.Range("D51").Value = "1.00"
$ws.Range("E51").Value = "  +0.04%  "
